# Users Controller: renamed 'Promote' to 'Roles'. Added POST action and cleaned up form
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("M0 - Account Mgmt")

# Row 6/7: status moved from "In progress" to "Done" (Good style), the
# "Need to add controller POST ..." notes are no longer needed now that the
# POST action + view cleanup are done.
$ws.Range("C6").Value = "Done"
$ws.Range("C6").Style = "Good"
$ws.Range("D6").Clear()

$ws.Range("C7").Value = "Done"
$ws.Range("C7").Style = "Good"
$ws.Range("D7").Clear()

# Row 30 (Obscure/encrypt connection string note): clarify the real risk.
$ws.Range("D30").Value = "Connection string is in Web.Config --> Will show up on github. "

# Row 18 (View Users List note): also mention filtering by game.
$ws.Range("D18").Value = "Should to filter by game (Warmachine by default). Also think about how to store/retrieve score summary data. "

# Move the cursor/selection to D18, matching the saved view state.
$ws.Range("D18").Select()
